# Update NATMI LR-pair output with refreshed TPM-derived statistics
# (Inha -> Tgfbr3), per "update scripts wuth new tpm" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> ECs (receptor average/total expr & derived specificities, edge weights)
$ws.Range("M2").Value = 24.412876
$ws.Range("N2").Value = 73.23862800000001
$ws.Range("O2").Value = 0.2909048961728503
$ws.Range("P2").Value = 0.2909048961728504
$ws.Range("Q2").Value = 0.8519279961466667
$ws.Range("R2").Value = 7.667351965320001
$ws.Range("S2").Value = 0.2909048961728503
$ws.Range("T2").Value = 0.2909048961728504

# Row 3: FAPs -> FAPs (only derived specificity columns shift, values re-normalized)
$ws.Range("O3").Value = 0.6669924168760825
$ws.Range("P3").Value = 0.6669924168760826
$ws.Range("S3").Value = 0.6669924168760825
$ws.Range("T3").Value = 0.6669924168760826

# Row 4: FAPs -> Inflammatory-Mac (receptor-expressing cell count drops 2 -> 1)
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.06139466666666666
$ws.Range("N4").Value = 0.184184
$ws.Range("O4").Value = 0.0007315815282162339
$ws.Range("P4").Value = 0.000731581528216234
$ws.Range("Q4").Value = 0.002142469217777777
$ws.Range("R4").Value = 0.01928222296
$ws.Range("S4").Value = 0.0007315815282162339
$ws.Range("T4").Value = 0.000731581528216234

# Row 5: FAPs -> MuSCs (receptor average/total expr & derived specificities, edge weights)
$ws.Range("M5").Value = 3.281552666666667
$ws.Range("N5").Value = 9.844658000000001
$ws.Range("O5").Value = 0.03910312483389531
$ws.Range("P5").Value = 0.03910312483389532
$ws.Range("Q5").Value = 0.1145152495577778
$ws.Range("R5").Value = 1.03063724602
$ws.Range("S5").Value = 0.03910312483389531
$ws.Range("T5").Value = 0.03910312483389532

# Row 6: FAPs -> Neutrophils (receptor average/total expr & derived specificities, edge weights)
$ws.Range("M6").Value = 0.1821903333333333
$ws.Range("N6").Value = 0.546571
$ws.Range("O6").Value = 0.002170987965614143
$ws.Range("P6").Value = 0.002170987965614143
$ws.Range("Q6").Value = 0.006357835332222223
$ws.Range("R6").Value = 0.05722051799
$ws.Range("S6").Value = 0.002170987965614143
$ws.Range("T6").Value = 0.002170987965614143

# Row 7: FAPs -> Resolving-Mac (receptor-expressing cell count drops 3 -> 1)
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.008139666666666667
$ws.Range("N7").Value = 0.024419
$ws.Range("O7").Value = 0.00009699262334139891
$ws.Range("P7").Value = 0.00009699262334139893
$ws.Range("Q7").Value = 0.0002840472344444444
$ws.Range("R7").Value = 0.00255642511
$ws.Range("S7").Value = 0.00009699262334139891
$ws.Range("T7").Value = 0.00009699262334139893
